# Updated the construction and residential pages buy buttons.
#
# This adds a "qty" / computed "price" buy-button block (columns G:H) to the
# "conscruction_services" worksheet, and moves the active/selected tab from
# "Aerial Photography Worksheet" to "conscruction_services".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conscruction_services")

# ---------------------------------------------------------------------
# Header row (G1:H1) - two new shared strings
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "qty"
$ws.Range("H1").Value = "price (Built: Up to 300K SF | Construction: Up to 30 Acres)"

# ---------------------------------------------------------------------
# Quantity column G2:G23 -> 1..22
#   G2 is a literal 1, G3 is a literal SUM formula, G4:G39 is a (shared)
#   "add one to the cell above" formula.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = 1
$ws.Range("G3").Formula = "=SUM(G2,1)"
for ($r = 4; $r -le 23; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 7).Formula = "=SUM(G$prev,1)"
}

# ---------------------------------------------------------------------
# Price column H2:H23 -> qty * unit price (B2 for qty 1-4, C2 for qty
# 5-10, D2 for qty 11-22), currency formatted + left aligned.
# ---------------------------------------------------------------------
$priceRefs = @{
    2  = "B2"; 3  = "B2"; 4  = "B2"; 5  = "B2"
    6  = "C2"; 7  = "C2"; 8  = "C2"; 9  = "C2"; 10 = "C2"; 11 = "C2"
    12 = "D2"; 13 = "D2"; 14 = "D2"; 15 = "D2"; 16 = "D2"; 17 = "D2"
    18 = "D2"; 19 = "D2"; 20 = "D2"; 21 = "D2"; 22 = "D2"; 23 = "D2"
}
for ($r = 2; $r -le 23; $r++) {
    $ref = $priceRefs[$r]
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = "=" + $ref + "*G" + $r
    $cell.NumberFormat = """$""#,##0.00"
    $cell.HorizontalAlignment = -4131
}

# ---------------------------------------------------------------------
# Trailing blank, currency-formatted cells H24:H39 (same "$" style as
# the empty B8:D10 block further up the sheet, no qty column alongside).
# ---------------------------------------------------------------------
for ($r = 24; $r -le 39; $r++) {
    $ws.Cells.Item($r, 8).NumberFormat = """$""#,##0.00"
}

# ---------------------------------------------------------------------
# Column widths for the new G/H columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 2.75
$ws.Columns.Item(8).ColumnWidth = 48.75

# ---------------------------------------------------------------------
# Move the active tab / selection from "Aerial Photography Worksheet"
# to "conscruction_services".
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
